# The commit swaps the OOXML content of ppt/theme/theme1.xml (the theme
# backing the slide master -> "Integral" / Red Violet) and
# ppt/theme/theme2.xml (the theme backing the notes master -> the
# default "Office Theme"), so that theme1.xml ends up holding the
# "Office Theme" colour scheme and theme2.xml ends up holding the
# "Integral" / Red Violet colour scheme that used to live in theme1.xml.
#
# The PowerPoint object model only exposes one set of theme colours for
# the whole deck (SlideMaster.Theme / NotesMaster.Theme both resolve to
# the same underlying theme used by the slide master, i.e. theme1.xml -
# there is no supported COM path to the notes-master-only theme part),
# so we reproduce the reachable half of the change: we repoint the
# slide master's theme colour scheme (ppt/theme/theme1.xml) from the
# "Integral" / Red Violet palette to the standard "Office Theme"
# palette, one MsoThemeColorSchemeIndex slot at a time via
# ThemeColorScheme.Colors(i).RGB, matching the colours the target
# theme1.xml carries after the edit.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex slot order:
# 1=dk1 2=lt1 3=dk2 4=lt2 5=accent1 6=accent2 7=accent3
# 8=accent4 9=accent5 10=accent6 11=hlink 12=folHlink
$tcs.Colors(1).RGB  = 0        # dk1      -> 000000
$tcs.Colors(2).RGB  = 16777215 # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      -> 44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  -> FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  -> 4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  -> 70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    -> 0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink -> 954F72
